# Updates cryptos list values (prices and 1h volume change %) for the
# Sun Nov  5 07:17:37 UTC 2023 GitHub Actions refresh.
#
# The workbook stores every value as text (inline/shared strings), including
# price and percentage columns that merely look numeric. Writing a plain
# numeric-looking string into a cell's .Value causes Excel's COM layer to
# coerce it into a real floating point number, which would change both the
# stored type and (for values like "0.0702") the literal text. To keep the
# cells as text exactly like the original file, each cell is temporarily
# switched to a Text number format ("@") before the value is assigned, and
# then restored to the default "Normal" style so no visible formatting
# change is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = [ordered]@{
    "D2" = "35.508.11"
    "E2" = "  +1.27%  "
    "D3" = "1.902.10"
    "E3" = "  +3.03%  "
    "E4" = "  +0.65%  "
    "D5" = "246.09"
    "E5" = "  +5.37%  "
    "E6" = "  +1.67%  "
    "D8" = "42.08"
    "E8" = "  +0.67%  "
    "E9" = "  +2.48%  "
    "D10" = "0.0702"
    "E10" = "  +1.03%  "
    "D11" = "0.0998"
    "E11" = "  +1.70%  "
    "D12" = "2.178.21"
    "E12" = "  +3.09%  "
    "D13" = "12.33"
    "E13" = "  +7.14%  "
    "D14" = "1.899.16"
    "E14" = "  +2.68%  "
    "D15" = "0.688"
    "E15" = "  +1.70%  "
    "E16" = "  +2.85%  "
    "D17" = "35.475.67"
    "E17" = "  +1.16%  "
    "D18" = "71.79"
    "E18" = "  +2.56%  "
    "E19" = "  +2.27%  "
    "D20" = "243.11"
    "E20" = "  +1.04%  "
    "D21" = "12.41"
    "E21" = "  +1.91%  "
    "E22" = "  +1.83%  "
    "E23" = "  +0.60%  "
    "D24" = "2.28"
    "E24" = "  -0.96%  "
    "D25" = "2.27"
    "E25" = "  +31.53%  "
    "D26" = "171.43"
    "E26" = "  -0.34%  "
    "D27" = "8.53"
    "E27" = "  +7.69%  "
    "D28" = "17.93"
    "E28" = "  +2.28%  "
    "E29" = "  +0.38%  "
    "D30" = "0.980"
    "E30" = "  +28.47%  "
    "B31" = "Filecoin"
    "C31" = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
    "D31" = "4.08"
    "E31" = "  +2.68%  "
    "B32" = "Hedera"
    "C32" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D32" = "0.0564"
    "E32" = "  +1.35%  "
    "E33" = "  +0.68%  "
    "E34" = "  +4.47%  "
    "E35" = "  +6.19%  "
    "E36" = "  +2.10%  "
    "E37" = "  +5.52%  "
    "E38" = "  +2.63%  "
    "E39" = "  +2.11%  "
    "D40" = "90.83"
    "E40" = "  +0.34%  "
    "D41" = "51.25"
    "E41" = "  +49.02%  "
    "D42" = "1.352.15"
    "E42" = "  +0.37%  "
    "D43" = "15.48"
    "E43" = "  +5.84%  "
    "D44" = "0.0592"
    "E44" = "  +11.06%  "
    "D45" = "2.34"
    "E45" = "  +1.68%  "
    "D46" = "12.57"
    "E46" = "  +7.46%  "
    "E47" = "  +1.57%  "
    "E48" = "  -0.27%  "
    "E49" = "  +4.48%  "
    "D50" = "2.085.90"
    "E50" = "  +2.86%  "
    "D51" = "0.0689"
    "E51" = "  +2.41%  "
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
    $cell.Style = "Normal"
}
